$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data values (row 2 and row 3) ---
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# --- Add a new (blank) row 6, matching the style/formatting of the data rows above ---
$ws.Range("A6:C6").WrapText = 1

# --- Re-fit the column widths to their (now shorter) numeric content ---
$ws.Columns("A:C").AutoFit()

# --- Move the active selection to C5, matching the saved cursor position ---
$ws.Range("C5").Select()
